# The caption paragraph under the matrix pictures originally consisted of
# three separate runs (split apart by two <w:proofErr/> "grammar" markers
# that LibreOffice/Word's grammar checker had inserted around the word
# "Matrix)"):
#
#   "(Original "  +  "Matrix)   "  +  "                                                    (Reduced Row Echelon Form matrix)"
#
# The edit collapses this into a single run/single <w:t> containing the
# same caption text (with the run of spaces between the two parenthetical
# phrases reduced by one character), and drops the now-stale proofErr
# markers. We reproduce that by searching for the old concatenated text
# and replacing it with the new single-run text.

$d = $word.ActiveDocument

# Build the exact old/new strings without relying on PowerShell's "*"
# string-repeat operator (not supported by this host - it coerces to a
# numeric multiply), using PadRight instead.
$oldText = "(Original Matrix)" + "".PadRight(55, ' ') + "(Reduced Row Echelon Form matrix)"
$newText = "(Original Matrix)" + "".PadRight(54, ' ') + "(Reduced Row Echelon Form matrix)"

$found = $d.Content.Find.Execute(
    $oldText,   # FindText
    $true,      # MatchCase
    $false,     # MatchWholeWord
    $false,     # MatchWildcards
    $false,     # MatchSoundsLike
    $false,     # MatchAllWordForms
    $true,      # Forward
    1,          # Wrap (wdFindContinue)
    $false,     # Format
    $newText,   # ReplaceWith
    2)          # Replace (wdReplaceAll)

if (-not $found) {
    throw "Could not find the '(Original Matrix) ... (Reduced Row Echelon Form matrix)' caption text to replace."
}
